$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "Action"
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item("Severity").Name = "Action"
$ws.Range("A1").Select()
